$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 6): new columns M, O, Q ---
$ws.Range("M6").Value = "precision"
$ws.Range("O6").Value = "log-based"
$ws.Range("Q6").Value = "trace fitness"

# --- Row 7: standalone (non-shared) formula in O7 ---
$ws.Range("M7").Value = 0.4736
$ws.Range("O7").Formula = "=M7*2"
$ws.Range("Q7").Value = 0.9986

# --- Rows 8-24: M (precision), O (shared formula =M*2), Q (trace fitness) ---
$ws.Range("M8").Value = 0.4791
$ws.Range("M9").Value = 0.4736
$ws.Range("M10").Value = 0.4615
$ws.Range("M11").Value = 0.4791
$ws.Range("M12").Value = 0
$ws.Range("M13").Value = 0.4615
$ws.Range("M14").Value = 0.4736
$ws.Range("M15").Value = 0.4791
$ws.Range("M16").Value = 0.4791
$ws.Range("M17").Value = 0.3571
$ws.Range("M18").Value = 0.3571
$ws.Range("M19").Value = 0.4791
$ws.Range("M20").Value = 0.1666
$ws.Range("M21").Value = 0.3823
$ws.Range("M22").Value = 0
$ws.Range("M23").Value = 0.4791
$ws.Range("M24").Value = 0.1666

# Shared formula group covering O8:O24 (written in one shot so the engine
# emits a single shared-formula group, matching si="0" ref="O8:O24")
$ws.Range("O8:O24").Formula = "=M8*2"

$ws.Range("Q8").Value = 0.9978
$ws.Range("Q9").Value = 0.9986
$ws.Range("Q10").Value = 0.9984
$ws.Range("Q11").Value = 0.9995
$ws.Range("Q12").Value = 1
$ws.Range("Q13").Value = 0.9987
$ws.Range("Q14").Value = 0.9986
$ws.Range("Q15").Value = 0.9995
$ws.Range("Q16").Value = 0.9981
$ws.Range("Q17").Value = 0.9989
$ws.Range("Q18").Value = 0.9994
$ws.Range("Q19").Value = 0.9981
$ws.Range("Q20").Value = 1
$ws.Range("Q21").Value = 0.99907
$ws.Range("Q22").Value = 1
$ws.Range("Q23").Value = 0.9979
$ws.Range("Q24").Value = 1

# --- View: selection moves to P17, and the prior scroll-anchor (topLeftCell)
#     is cleared by simply re-selecting within the now-visible range. ---
[void]$ws.Range("P17").Select()
